# Update "Mispronounced words" note: add 5 new word entries (rows 8-14)
# specific / niche / bargaining / suite / automate / modular / character

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Pronunciation {
    param($Cell, $Middle)
    $full = "/" + $Middle + "/"
    $Cell.Value = $full
    $midLen = $Middle.Length
    # characters 2..(1+midLen) get the GWIPA phonetic font
    $Cell.Characters(2, $midLen).Font.Name = "GWIPA"
    # trailing slash explicitly kept in the default (Calibri) font
    $Cell.Characters(2 + $midLen, 1).Font.Name = "Calibri"
}

# Row 8: specific
$ws.Range("A8").Value = "specific"
Set-Pronunciation $ws.Range("B8") "spE'sIfIk"

# Row 9: niche
$ws.Range("A9").Value = "niche"
Set-Pronunciation $ws.Range("B9") "ni:S"

# Row 10: bargaining [U.]
$ws.Range("A10").Value = "bargaining"
Set-Pronunciation $ws.Range("B10") "'bA:gEnIN"
$ws.Range("C10").Value = "[U.]"

# Row 11: suite / n.  (boxed word pair with suite above automate)
$b11 = $ws.Range("B11")
$b11.Value = " /swi:t/ "
$b11.Font.Name = "GWIPA"

$blank = $ws.Range("Z50")
$blank.Copy()
$a11 = $ws.Range("A11")
$a11.PasteSpecial(-4122)   # xlPasteFormats - start from a clean, unformatted cell
$a11.Value = "suite"
$a11.Font.Name = "Arial"
$a11.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$a11.Borders.Item(8).LineStyle = 1    # xlEdgeTop

$ws.Range("C11").Value = "n."

# Row 12: automate / vt.
$b12 = $ws.Range("B12")
$b12.Value = " /'O:tEmeIt/ "
$b12.Font.Name = "GWIPA"

$blank.Copy()
$a12 = $ws.Range("A12")
$a12.PasteSpecial(-4122)
$a12.Value = "automate"
$a12.Font.Name = "Arial"
$a12.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$a12.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$a12.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$ws.Range("C12").Value = "vt."

# Row 13: modular / adj.
$ws.Range("A13").Value = "modular"
Set-Pronunciation $ws.Range("B13") "'mOdjulE"
$ws.Range("C13").Value = "adj."

# Row 14: character / n.single
$ws.Range("A14").Value = "character"
Set-Pronunciation $ws.Range("B14") "'kQrEktE"
$ws.Range("C14").Value = "n.single"

# Column C width tweak + selection cursor (cosmetic, matches the saved-file state)
$ws.Columns("C").ColumnWidth = 14.28515625
$ws.Range("D13").Select()

Write-Host "Applied Mispronounced words update (rows 8-14)."
